$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bvals = New-Object "object[,]" 24,1
$bvals[0,0] = 1.199213271610233
$bvals[1,0] = 1.049321093537515
$bvals[2,0] = 0.9569891392220029
$bvals[3,0] = 0.9192903332054243
$bvals[4,0] = 0.9130261397797312
$bvals[5,0] = 0.9564810119189815
$bvals[6,0] = 1.14759337227639
$bvals[7,0] = 1.519932936973134
$bvals[8,0] = 1.791942598005846
$bvals[9,0] = 1.915338907632815
$bvals[10,0] = 1.962015122589662
$bvals[11,0] = 1.951964878313788
$bvals[12,0] = 1.91918002760093
$bvals[13,0] = 1.899091619540798
$bvals[14,0] = 1.783871293892219
$bvals[15,0] = 1.713098210472424
$bvals[16,0] = 1.672359312328695
$bvals[17,0] = 1.658560373303601
$bvals[18,0] = 1.720635462592952
$bvals[19,0] = 1.928811141306539
$bvals[20,0] = 2.064566090058918
$bvals[21,0] = 1.992139235232116
$bvals[22,0] = 1.717228028083412
$bvals[23,0] = 1.419472465225738
$ws.Range("B2:B25").Value2 = $bvals

$cvals = New-Object "object[,]" 24,1
$cvals[0,0] = 0.354944670868008
$cvals[1,0] = 0.3110170542605033
$cvals[2,0] = 0.2839267843903315
$cvals[3,0] = 0.2728579909621658
$cvals[4,0] = 0.2710182727439019
$cvals[5,0] = 0.2837776246635428
$cvals[6,0] = 0.3398233259204062
$cvals[7,0] = 0.4487710321283203
$cvals[8,0] = 0.5282161006501269
$cvals[9,0] = 0.564225176901175
$cvals[10,0] = 0.5778416822854524
$cvals[11,0] = 0.574909993150527
$cvals[12,0] = 0.5653458057975058
$cvals[13,0] = 0.5594849293671018
$cvals[14,0] = 0.5258601436306094
$cvals[15,0] = 0.5051985243355261
$cvals[16,0] = 0.4933022189781013
$cvals[17,0] = 0.4892722401429523
$cvals[18,0] = 0.5073992646434249
$cvals[19,0] = 0.5681555687934292
$cvals[20,0] = 0.6077502989650725
$cvals[21,0] = 0.5866283605292324
$cvals[22,0] = 0.5064043649381915
$cvals[23,0] = 0.4194019072779156
$ws.Range("C2:C25").Value2 = $cvals

$dvals = New-Object "object[,]" 24,1
$dvals[0,0] = 0.01503616514069961
$dvals[1,0] = 0.01326176328310424
$dvals[2,0] = 0.01216786834112327
$dvals[3,0] = 0.01172102449260137
$dvals[4,0] = 0.01164676265852904
$dvals[5,0] = 0.0121618463455988
$dvals[6,0] = 0.0144252858462508
$dvals[7,0] = 0.01882763574549529
$dvals[8,0] = 0.02203855149714684
$dvals[9,0] = 0.02349391636342801
$dvals[10,0] = 0.02404423550894563
$dvals[11,0] = 0.02392575037173827
$dvals[12,0] = 0.02353920759617267
$dvals[13,0] = 0.02330233432530093
$dvals[14,0] = 0.02194333044912611
$dvals[15,0] = 0.02110824433166414
$dvals[16,0] = 0.0206274283238983
$dvals[17,0] = 0.02046454805777387
$dvals[18,0] = 0.02119719231695427
$dvals[19,0] = 0.02365276645058856
$dvals[20,0] = 0.02525297134520343
$dvals[21,0] = 0.02439934891990703
$dvals[22,0] = 0.02115698115353126
$dvals[23,0] = 0.01764071084648577
$ws.Range("D2:D25").Value2 = $dvals

$evals = New-Object "object[,]" 24,1
$evals[0,0] = 0.4254321993331303
$evals[1,0] = 0.3710241705299779
$evals[2,0] = 0.3377147263827851
$evals[3,0] = 0.3241633335824616
$evals[4,0] = 0.3219144347038565
$evals[5,0] = 0.3375318791631088
$evals[6,0] = 0.4066508972364602
$evals[7,0] = 0.5430615204112854
$evals[8,0] = 0.6439546468780861
$evals[9,0] = 0.690031129686048
$evals[10,0] = 0.7075073357755883
$evals[11,0] = 0.7037422395795971
$evals[12,0] = 0.6914683329632823
$evals[13,0] = 0.6839539284820262
$evals[14,0] = 0.6409472439179638
$evals[15,0] = 0.6146116184445276
$evals[16,0] = 0.5994808435457628
$evals[17,0] = 0.5943606385196745
$evals[18,0] = 0.617413338128884
$evals[19,0] = 0.695072698124676
$evals[20,0] = 0.7459920709235348
$evals[21,0] = 0.7187996576595594
$evals[22,0] = 0.6161466499250565
$evals[23,0] = 0.506051184457263
$ws.Range("E2:E25").Value2 = $evals

$fvals = New-Object "object[,]" 24,1
$fvals[0,0] = 0.4030827437275448
$fvals[1,0] = 0.3993247100541879
$fvals[2,0] = 0.3975013686066049
$fvals[3,0] = 0.3968793584726313
$fvals[4,0] = 0.3967833622882466
$fvals[5,0] = 0.3974924909217421
$fvals[6,0] = 0.4016860639633606
$fvals[7,0] = 0.4137842399471268
$fvals[8,0] = 0.4250827099878194
$fvals[9,0] = 0.4307563817382345
$fvals[10,0] = 0.4329824448454502
$fvals[11,0] = 0.4324995599823183
$fvals[12,0] = 0.4309379622260394
$fvals[13,0] = 0.4299915633601756
$fvals[14,0] = 0.4247227336731996
$fvals[15,0] = 0.4216277917657152
$fvals[16,0] = 0.4198978730123741
$fvals[17,0] = 0.4193207519033493
$fvals[18,0] = 0.4219520513896455
$fvals[19,0] = 0.4313945299252424
$fvals[20,0] = 0.4380182199325873
$fvals[21,0] = 0.434441369289047
$fvals[22,0] = 0.4218052997839123
$fvals[23,0] = 0.4100916148856655
$ws.Range("F2:F25").Value2 = $fvals

$hvals = New-Object "object[,]" 24,1
$hvals[0,0] = 0.07973214163530429
$hvals[1,0] = 0.07973214163530429
$hvals[2,0] = 0.07973214163530429
$hvals[3,0] = 0.07973214163530429
$hvals[4,0] = 0.07973214163530429
$hvals[5,0] = 0.07973214163530429
$hvals[6,0] = 0.07973214163530429
$hvals[7,0] = 0.07973214163530429
$hvals[8,0] = 0.07973214163530429
$hvals[9,0] = 0.07973214163530429
$hvals[10,0] = 0.07973214163530429
$hvals[11,0] = 0.07973214163530429
$hvals[12,0] = 0.07973214163530429
$hvals[13,0] = 0.07973214163530429
$hvals[14,0] = 0.07973214163530429
$hvals[15,0] = 0.07973214163530429
$hvals[16,0] = 0.07973214163530429
$hvals[17,0] = 0.07973214163530429
$hvals[18,0] = 0.07973214163530429
$hvals[19,0] = 0.07973214163530429
$hvals[20,0] = 0.07973214163530429
$hvals[21,0] = 0.07973214163530429
$hvals[22,0] = 0.07973214163530429
$hvals[23,0] = 0.07973214163530429
$ws.Range("H2:H25").Value2 = $hvals

$ivals = New-Object "object[,]" 24,1
$ivals[0,0] = 0.2470587248881824
$ivals[1,0] = 0.2546024141071861
$ivals[2,0] = 0.2595819383436453
$ivals[3,0] = 0.2616981020691878
$ivals[4,0] = 0.2620547286666621
$ivals[5,0] = 0.2596101261489689
$ivals[6,0] = 0.2495873564182904
$ivals[7,0] = 0.2327115326399998
$ivals[8,0] = 0.2220342732705323
$ivals[9,0] = 0.2175565732668723
$ivals[10,0] = 0.2159160425858815
$ivals[11,0] = 0.2162669024471313
$ivals[12,0] = 0.2174204990424968
$ivals[13,0] = 0.2181342982195957
$ivals[14,0] = 0.2223345807100117
$ivals[15,0] = 0.2250088668850321
$ivals[16,0] = 0.2265827293310387
$ivals[17,0] = 0.2271217251800426
$ivals[18,0] = 0.2247204877097619
$ivals[19,0] = 0.2170801608760939
$ivals[20,0] = 0.2124080368787808
$ivals[21,0] = 0.2148720707879548
$ivals[22,0] = 0.2248507506545945
$ivals[23,0] = 0.2369765749567776
$ws.Range("I2:I25").Value2 = $ivals

$ovals = New-Object "object[,]" 24,1
$ovals[0,0] = 1.244608197349265
$ovals[1,0] = 1.246152473051836
$ovals[2,0] = 1.248714745250382
$ovals[3,0] = 1.25016221876777
$ovals[4,0] = 1.250426856040676
$ovals[5,0] = 1.248732636715403
$ovals[6,0] = 1.244804124213232
$ovals[7,0] = 1.250024152912715
$ovals[8,0] = 1.26190253960533
$ovals[9,0] = 1.269088469333695
$ovals[10,0] = 1.272068743956197
$ovals[11,0] = 1.271415319235615
$ovals[12,0] = 1.269328449053091
$ovals[13,0] = 1.268084009353515
$ovals[14,0] = 1.261469028813877
$ovals[15,0] = 1.257869465458867
$ovals[16,0] = 1.255966690492045
$ovals[17,0] = 1.255351143183873
$ovals[18,0] = 1.258235278979754
$ovals[19,0] = 1.269934358543964
$ovals[20,0] = 1.279092037482741
$ovals[21,0] = 1.27406516419245
$ovals[22,0] = 1.258069375765729
$ovals[23,0] = 1.247211425952941
$ws.Range("O2:O25").Value2 = $ovals

